$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric/percent-looking values in columns D and E are written as literal text
# (matching the workbook's existing inlineStr storage), by pre-formatting those cells as Text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '296.97'
$ws.Range("E2").Value = '-3.70%'
$ws.Range("D3").Value = '40.68'
$ws.Range("E3").Value = '-0.92%'
$ws.Range("D4").Value = '5.038'
$ws.Range("E4").Value = '-3.29%'
$ws.Range("D5").Value = '0.07432'
$ws.Range("E5").Value = '-3.23%'
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").Value = '4.312'
$ws.Range("E6").Value = '0.38%'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = '1.580'
$ws.Range("E7").Value = '-3.88%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '0.9251'
$ws.Range("E8").Value = '1.13%'
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").Value = '2.397'
$ws.Range("E9").Value = '-1.36%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '0.1164'
$ws.Range("E10").Value = '-5.19%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1748'
$ws.Range("E11").Value = '-4.22%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.08770'
$ws.Range("E12").Value = '-4.55%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.04192'
$ws.Range("E13").Value = '-1.56%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.1053'
$ws.Range("E14").Value = '0.24%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001265'
$ws.Range("E15").Value = '0.67%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '0.006005'
$ws.Range("E16").Value = '2.98%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.361'
$ws.Range("E17").Value = '0.58%'
$ws.Range("E18").Value = '0.36%'
$ws.Range("D19").Value = '7.627'
$ws.Range("E19").Value = '4.17%'
$ws.Range("D20").Value = '0.1358'
$ws.Range("E20").Value = '-1.72%'
$ws.Range("D22").Value = '0.03860'
$ws.Range("E22").Value = '-4.08%'
$ws.Range("D23").Value = '0.001290'
$ws.Range("E23").Value = '2.14%'
$ws.Range("D24").Value = '0.003631'
$ws.Range("E24").Value = '-17.00%'
$ws.Range("D25").Value = '0.0001307'
$ws.Range("E25").Value = '0.49%'
$ws.Range("D26").Value = '0.0003744'
$ws.Range("D38").Value = '0.02314'
$ws.Range("D39").Value = '0.05013'
$ws.Range("E39").Value = '-5.69%'
$ws.Range("D40").Value = '0.007741'
$ws.Range("E40").Value = '-1.33%'
$ws.Range("E41").Value = '124.04%'
$ws.Range("D42").Value = '0.1278'
$ws.Range("E42").Value = '-2.61%'
$ws.Range("D43").Value = '0.007425'
$ws.Range("E43").Value = '11.51%'
$ws.Range("D44").Value = '0.007828'
$ws.Range("E44").Value = '-2.26%'
$ws.Range("D45").Value = '0.3221'
$ws.Range("E45").Value = '5.08%'
$ws.Range("D46").Value = '0.00006490'
$ws.Range("E46").Value = '-3.51%'
$ws.Range("D47").Value = '0.00000000754'
$ws.Range("E47").Value = '0.51%'
$ws.Range("E48").Value = '-11.09%'
$ws.Range("D49").Value = '0.004224'
$ws.Range("E49").Value = '36.16%'
$ws.Range("D50").Value = '0.00002112'
$ws.Range("E50").Value = '0.51%'
$ws.Range("D51").Value = '0.0002012'
$ws.Range("E51").Value = '0.51%'
